$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name text on both sheets
$wsInput.Range("B1").Value = "4132-RBI-EPP-DB-DL-NOREC-MOREREPAY-1st"
$wsOutput.Range("B1").Value = "4132-RBI-EPP-DB-DL-NOREC-MOREREPAY-1st"

# Update shortname on input sheet to a text value
$wsInput.Range("B2").Value = "413r"

# Update selection on input sheet
$wsInput.Range("B3").Select()

# Update selection on output sheet (stays B1)
$wsOutput.Range("B1").Select()

# Make output sheet the active sheet/tab
$wsOutput.Activate()
